# "updated naming conventions of database"
#
# The workbook models a small relational schema. Every table/column that
# was named after the "Crisis" entity is renamed to "Event" (crisisID ->
# eventID, Crisis -> Event, CrisisPerson -> EventPerson, CrisisOrganization
# -> EventOrganization, CitationCrisis -> CitationEvent, LinkCrisis ->
# LinkEvent, EmbedCrisis -> EmbedEvent). No other cell content, formatting,
# or layout changes.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# --- Rename "Crisis" entity/table/column names to "Event" -------------------
# "Crisis" section heading (A1) -> "Event"
$ws.Range("A1").Value = "Event"

# "CrisisPerson" join-table heading (A25) -> "EventPerson"
$ws.Range("A25").Value = "EventPerson"

# "CrisisOrganization" join-table heading (A29) -> "EventOrganization"
$ws.Range("A29").Value = "EventOrganization"

# "crisisID" foreign/primary key column, every place it appears
$ws.Range("D7").Value = "eventID"
$ws.Range("G9").Value = "eventID"
$ws.Range("D24").Value = "eventID"
$ws.Range("A26").Value = "eventID"
$ws.Range("A30").Value = "eventID"

# "CitationCrisis" join table name -> "CitationEvent"
$ws.Range("D5").Value = "CitationEvent"

# "LinkCrisis" join table name -> "LinkEvent"
$ws.Range("D22").Value = "LinkEvent"

# "EmbedCrisis" join table name -> "EmbedEvent"
$ws.Range("G7").Value = "EmbedEvent"

# --- Reset the saved cursor position back to the top-left cell --------------
# (the authored file no longer pins the selection to D33)
$null = $ws.Range("A1").Select()
